$wb = $excel.ActiveWorkbook

# Sheet 1: 01_리그테이블 — delete the "코셈" listing row (row 15)
$ws1 = $wb.Worksheets.Item("01_리그테이블")
$ws1.Rows.Item(15).Delete()

# Sheet 2: 02_통합집계_Rawdata — delete the "코셈" listing row (row 14)
$ws2 = $wb.Worksheets.Item("02_통합집계_Rawdata")
$ws2.Rows.Item(14).Delete()

# Sheet 3: 03_IPO현황_Summary — delete the "코셈" / "키움" listing row (row 12)
$ws3 = $wb.Worksheets.Item("03_IPO현황_Summary")
$ws3.Rows.Item(12).Delete()
